$wb = $excel.ActiveWorkbook

# ----- Sheet 1: tweet/data pull sheet -----
$ws1 = $wb.Worksheets.Item("Sheet 1")

# Insert a new "News" column between keyPhrases (J) and Irrelevent (K),
# pushing Irrelevent out to column L.
$ws1.Columns.Item(11).Insert()
$ws1.Range("K1").Value = "News"
$ws1.Range("K2").Value = 0

# New data pull values for the single data row.
$ws1.Range("A2").Value = "2019-06-22 07:23:42"
$ws1.Range("B2").Value = "We are unstoppable!! #ceat #Bhandup employee #volunteers continue to work towards improving relationship with the https://t.co/1fDScnkF43"
$ws1.Range("C2").Value = 8
$ws1.Range("D2").Value = 0
$ws1.Range("E2").Value = 0.7872353792190552
$ws1.Range("F2").Value = "plastic"
$ws1.Range("G2").Value = "Sustainability"
$ws1.Range("J2").Value = "['volunteers', 'ceat', 'Bhandup employee', 'relationship']"
$ws1.Range("L2").Value = 0

# ----- Sheet 2: Category lookup list -----
$ws2 = $wb.Worksheets.Item("Sheet 2")

# Drop the "Lifestyle" category (row 2) and close the gap.
$ws2.Rows.Item(2).Delete()

# Renumber the remaining index column (A) back to a contiguous 0-based range.
$ws2.Range("A2").Value = 0
$ws2.Range("A3").Value = 1
$ws2.Range("A4").Value = 2
$ws2.Range("A5").Value = 3
$ws2.Range("A6").Value = 4
$ws2.Range("A7").Value = 5
$ws2.Range("A8").Value = 6
$ws2.Range("A9").Value = 7
